$d = $word.ActiveDocument

# The title paragraph currently reads:
#   "MINI" (run) + bookmarkStart/_GoBack + bookmarkEnd + "SUMO" (run)
# It needs to become:
#   "MINI" (run) + "SUMITO" (new, separate run) + bookmarkStart/_GoBack + bookmarkEnd
# i.e. a new "SUMITO" run is inserted right before the bookmark, and the old
# "SUMO" run (which used to sit right after the bookmark) is removed.

$oldWord = "SUMO"
$newWord = "SUMITO"

$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range

# Insert the new text right before the bookmark. At this point it merges with the
# preceding "MINI" run's text (becomes one run reading "MINISUMITO") because plain
# text insertion always coalesces with a touching, identically-formatted run.
$bmRange.InsertBefore($newWord)

# Re-fetch the bookmark: its Range has shifted to the new position, right after the
# text we just inserted.
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range

# Split "SUMITO" back out into its own run by re-applying explicit character
# formatting (matching what is already in effect) to only that portion of text.
# Toggling a formatting property forces the engine to materialize an explicit,
# separate <w:r> for the touched range instead of leaving it merged.
$newRunRange = $d.Range($bmRange.Start - $newWord.Length, $bmRange.Start)
$newRunRange.Font.Bold = $false
$newRunRange.Font.Bold = $true

# Remove the old "SUMO" run, which now sits right after the bookmark.
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$oldRun = $d.Range($bmRange.End, $bmRange.End + $oldWord.Length)
$oldRun.Delete()
